$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on column A (the "code" column), which is unaffected.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Swap the contents of columns C (codeforiati:group-code) and D (codeforiati:group-name)
# for every row, including the header, so that:
#   - column C now holds what used to be in column D
#   - column D now holds what used to be in column C
for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
